$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Nutrition label placeholder -> real salad labels ----------------------
# Shared-string table order matters: the engine appends/reuses strings in the
# order cells are written, so touch F2 first (keeps "placeholder" -> "Cobb_Salad"
# at the same shared-string index), then F3, then F4 so the two brand new
# strings ("Caesar_Salad", "House_Salad") are appended right after it.
$ws.Range("F2").Value = "Cobb_Salad"
$ws.Range("F3").Value = "Caesar_Salad"
$ws.Range("F4").Value = "House_Salad"

# --- Column layout: narrow the Allergens / LocalIngredients columns --------
# (the host's ColumnWidth->stored-width conversion snaps to 1/6-character
# steps, so these land on the closest achievable width to the target
# 17.28515625 / 17.140625 from the source file)
$ws.Columns.Item(3).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 16.33
